$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.095493383854034164
$ws.Range("A2").Value = -0.0059999999672299964
$ws.Range("A3").Value = -0.0039999999749298354
$ws.Range("A4").Value = -0.0079999999527302634
$ws.Range("A5").Value = -0.0029999999787824194
$ws.Range("A6").Value = 0.054263324164343985
$ws.Range("A7").Value = -0.0099999999385858018
$ws.Range("A8").Value = -0.009999999935868864
$ws.Range("A9").Value = -0.0019999999748119279
$ws.Range("A10").Value = -0.0019999999710478278
$ws.Range("A11").Value = -0.0029999999652323694
$ws.Range("A12").Value = -0.0034999999617966182
$ws.Range("A13").Value = -0.0034999999580298535
$ws.Range("A14").Value = -0.007999999932845725
$ws.Range("A15").Value = -0.00099999996978628047
$ws.Range("A16").Value = -0.0019999999638051769
$ws.Range("A17").Value = -0.0019999999628694809
$ws.Range("A18").Value = -0.0039999999519269025
$ws.Range("A19").Value = -0.0039999999778115303
$ws.Range("A20").Value = -0.0039999999759405824
$ws.Range("A21").Value = -0.003999999975679458
$ws.Range("A22").Value = -0.0039999999754920523
$ws.Range("A23").Value = -0.0049999999677234896
$ws.Range("A24").Value = -0.019999999883599351
$ws.Range("A25").Value = -0.019999999881900266
$ws.Range("A26").Value = -0.0024999999798573924
$ws.Range("A27").Value = -0.0024999999796233574
$ws.Range("A28").Value = -0.0019999999813506975
$ws.Range("A29").Value = -0.0069999999541110469
$ws.Range("A30").Value = 0.000019425379887305638
$ws.Range("A31").Value = 0.017615383393732031
$ws.Range("A32").Value = -0.0099999999380564475
$ws.Range("A33").Value = -0.0039999999701443301
